# Update team-specific matrix percentages with figures recalculated
# from games pulled March 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.176
$ws.Range("C2").Value = 0.5866666666666667
$ws.Range("J2").Value = 0.01333333333333333
$ws.Range("P2").Value = 0.1173333333333333
$ws.Range("S2").Value = 0.1066666666666667
$ws.Range("B3").Value = 0.004484304932735426
$ws.Range("C3").Value = 0.008968609865470852
$ws.Range("J3").Value = 0.01345291479820628
$ws.Range("P3").Value = 0.7488789237668162
$ws.Range("S3").Value = 0.2242152466367713
$ws.Range("J4").Value = 0.03703703703703703
$ws.Range("P4").Value = 0.7222222222222222
$ws.Range("S4").Value = 0.2407407407407407
$ws.Range("B6").Value = 0.04477611940298507
$ws.Range("D6").Value = 0.004975124378109453
$ws.Range("F6").Value = 0.03482587064676617
$ws.Range("J6").Value = 0.2238805970149254
$ws.Range("O6").Value = 0.01492537313432836
$ws.Range("Q6").Value = 0.1940298507462687
$ws.Range("R6").Value = 0.07960199004975124
$ws.Range("S6").Value = 0.4029850746268657
$ws.Range("B7").Value = 0.1161825726141079
$ws.Range("D7").Value = 0.01244813278008299
$ws.Range("E7").Value = 0.004149377593360996
$ws.Range("F7").Value = 0.03319502074688797
$ws.Range("J7").Value = 0.1659751037344398
$ws.Range("O7").Value = 0.02489626556016597
$ws.Range("Q7").Value = 0.1618257261410788
$ws.Range("R7").Value = 0.07468879668049792
$ws.Range("S7").Value = 0.4066390041493776
$ws.Range("B8").Value = 0.1148775894538606
$ws.Range("D8").Value = 0.02824858757062147
$ws.Range("E8").Value = 0.001883239171374765
$ws.Range("F8").Value = 0.05649717514124294
$ws.Range("J8").Value = 0.1073446327683616
$ws.Range("O8").Value = 0.01129943502824859
$ws.Range("Q8").Value = 0.1789077212806026
$ws.Range("R8").Value = 0.09416195856873823
$ws.Range("S8").Value = 0.4067796610169492
$ws.Range("D9").Value = 0.01
$ws.Range("F9").Value = 0.045
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.155
$ws.Range("R9").Value = 0.125
$ws.Range("S9").Value = 0.44
$ws.Range("B10").Value = 0.1308937368050669
$ws.Range("D10").Value = 0.02251935256861365
$ws.Range("F10").Value = 0.06403940886699508
$ws.Range("J10").Value = 0.1048557353976073
$ws.Range("O10").Value = 0.01970443349753695
$ws.Range("Q10").Value = 0.1724137931034483
$ws.Range("R10").Value = 0.09992962702322308
$ws.Range("S10").Value = 0.3856439127375088
$ws.Range("G11").Value = 0.1425178147268409
$ws.Range("J11").Value = 0.1187648456057007
$ws.Range("K11").Value = 0.2256532066508314
$ws.Range("L11").Value = 0.498812351543943
$ws.Range("S11").Value = 0.01425178147268409
$ws.Range("G12").Value = 0.6863636363636364
$ws.Range("J12").Value = 0.2454545454545455
$ws.Range("K12").Value = 0.00909090909090909
$ws.Range("L12").Value = 0.02727272727272727
$ws.Range("S12").Value = 0.03181818181818181
$ws.Range("G13").Value = 0.6724137931034483
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("S13").Value = 0.05172413793103448
$ws.Range("F15").Value = 0.02109704641350211
$ws.Range("H15").Value = 0.1392405063291139
$ws.Range("I15").Value = 0.05907172995780591
$ws.Range("J15").Value = 0.3586497890295359
$ws.Range("K15").Value = 0.05485232067510549
$ws.Range("M15").Value = 0.01265822784810127
$ws.Range("O15").Value = 0.0379746835443038
$ws.Range("S15").Value = 0.3164556962025317
$ws.Range("H16").Value = 0.1844262295081967
$ws.Range("I16").Value = 0.0778688524590164
$ws.Range("J16").Value = 0.3565573770491803
$ws.Range("K16").Value = 0.1352459016393443
$ws.Range("M16").Value = 0.03688524590163934
$ws.Range("O16").Value = 0.02868852459016394
$ws.Range("S16").Value = 0.180327868852459
$ws.Range("F17").Value = 0.01126126126126126
$ws.Range("H17").Value = 0.1711711711711712
$ws.Range("I17").Value = 0.08108108108108109
$ws.Range("J17").Value = 0.4279279279279279
$ws.Range("K17").Value = 0.09684684684684684
$ws.Range("M17").Value = 0.01351351351351351
$ws.Range("O17").Value = 0.06081081081081081
$ws.Range("S17").Value = 0.1373873873873874
$ws.Range("F18").Value = 0.016
$ws.Range("H18").Value = 0.232
$ws.Range("I18").Value = 0.112
$ws.Range("J18").Value = 0.364
$ws.Range("K18").Value = 0.07199999999999999
$ws.Range("M18").Value = 0.016
$ws.Range("O18").Value = 0.068
$ws.Range("S18").Value = 0.12
$ws.Range("F19").Value = 0.00975292587776333
$ws.Range("H19").Value = 0.2106631989596879
$ws.Range("I19").Value = 0.06827048114434331
$ws.Range("J19").Value = 0.3556566970091027
$ws.Range("K19").Value = 0.1404421326397919
$ws.Range("M19").Value = 0.02535760728218465
$ws.Range("O19").Value = 0.0611183355006502
$ws.Range("S19").Value = 0.1287386215864759
